$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN": append 4 new rows (41-44) after the existing last row (40).
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Dates/sprint-name columns (A, C) must stay literal text, not get
# auto-converted to Excel date serials, so force text format before typing.
$wsAmsin.Range("A41:A44").NumberFormat = "@"
$wsAmsin.Range("C41:C44").NumberFormat = "@"

$wsAmsin.Range("A41").Value = "2022-08-03"
$wsAmsin.Range("B41").Value = 44776.69670299769
$wsAmsin.Range("C41").Value = "165_scndcycle"
$wsAmsin.Range("D41").Value = 119
$wsAmsin.Range("E41").Value = 119
$wsAmsin.Range("F41").Value = 0
$wsAmsin.Range("G41").Value = 2.66

$wsAmsin.Range("A42").Value = "2022-08-04"
$wsAmsin.Range("B42").Value = 44777.39890631945
$wsAmsin.Range("C42").Value = "165_finalrun"
$wsAmsin.Range("D42").Value = 119
$wsAmsin.Range("E42").Value = 119
$wsAmsin.Range("F42").Value = 0
$wsAmsin.Range("G42").Value = 2.56

$wsAmsin.Range("A43").Value = "2022-08-22"
$wsAmsin.Range("B43").Value = 44795.67989542824
$wsAmsin.Range("C43").Value = "166fstcycle"
$wsAmsin.Range("D43").Value = 119
$wsAmsin.Range("E43").Value = 119
$wsAmsin.Range("F43").Value = 0
$wsAmsin.Range("G43").Value = 2.64

$wsAmsin.Range("A44").Value = "2022-08-23"
$wsAmsin.Range("B44").Value = 44796.91313158565
$wsAmsin.Range("C44").Value = "166cyclescnd"
$wsAmsin.Range("D44").Value = 119
$wsAmsin.Range("E44").Value = 119
$wsAmsin.Range("F44").Value = 0
$wsAmsin.Range("G44").Value = 2.58

# Match the look of the preceding rows: copy the formatting of row 40 (which
# carries the sheet's normal style + the date/time display on column B) onto
# the freshly written rows.
$wsAmsin.Range("A40:G40").Copy()
$wsAmsin.Range("A41:G44").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "BETA": append 2 new rows (20-21) after the existing last row (19).
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

$wsBeta.Range("A20:A21").NumberFormat = "@"
$wsBeta.Range("C20:C21").NumberFormat = "@"

$wsBeta.Range("A20").Value = "2022-08-04"
$wsBeta.Range("B20").Value = 44777.572070625
$wsBeta.Range("C20").Value = "165beta"
$wsBeta.Range("D20").Value = 119
$wsBeta.Range("E20").Value = 119
$wsBeta.Range("F20").Value = 0
$wsBeta.Range("G20").Value = 2.8

$wsBeta.Range("A21").Value = "2022-08-24"
$wsBeta.Range("B21").Value = 44797.55309025463
$wsBeta.Range("C21").Value = "166_beta"
$wsBeta.Range("D21").Value = 119
$wsBeta.Range("E21").Value = 119
$wsBeta.Range("F21").Value = 0
$wsBeta.Range("G21").Value = 2.64

$wsBeta.Range("A19:G19").Copy()
$wsBeta.Range("A20:G21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "AMS": row 20 picks up the sheet's normal formatting (it previously
# had none) and its run-time value is corrected; rows 21-22 are then
# appended, with row 22 kept in the same un-styled state row 20 used to be
# in (only its Run Time cell carries the date/time style).
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Row 22 should end up looking the way row 20 looks *right now*, before any
# of today's edits (plain/un-styled, Run Time cell already date/time
# formatted). Stash that look on row 22 first (the clipboard keeps the copied
# formatting around even while other edits happen below); the values written
# into row 22 further down get stamped back with this same formatting at the
# very end, the same way row 20's original data must have been entered.
$wsAms.Range("A20:G20").Copy()
$wsAms.Range("A22:G22").PasteSpecial(-4122)

# Now bring row 20 itself in line with the rest of the sheet's data rows,
# and correct its Run Time value.
$wsAms.Range("A19:G19").Copy()
$wsAms.Range("A20").PasteSpecial(-4122)
$wsAms.Range("B20").Value = 44756.8288728125

# New row 21, styled like the rest of the sheet's data rows.
$wsAms.Range("A21").NumberFormat = "@"
$wsAms.Range("C21").NumberFormat = "@"
$wsAms.Range("A21").Value = "2022-08-04"
$wsAms.Range("B21").Value = 44777.82556201389
$wsAms.Range("C21").Value = "165_live"
$wsAms.Range("D21").Value = 119
$wsAms.Range("E21").Value = 119
$wsAms.Range("F21").Value = 0
$wsAms.Range("G21").Value = 2.6

$wsAms.Range("A19:G19").Copy()
$wsAms.Range("A21:G21").PasteSpecial(-4122)

# New row 22: write its values (forcing A/C to stay text so they don't turn
# into date serials), then restamp the plain look stashed above so the row
# ends up un-styled just like old row 20 was, while keeping the new values.
$wsAms.Range("A22").NumberFormat = "@"
$wsAms.Range("C22").NumberFormat = "@"
$wsAms.Range("A22").Value = "2022-08-24"
$wsAms.Range("B22").Value = 44797.92940078524
$wsAms.Range("C22").Value = "166_live"
$wsAms.Range("D22").Value = 119
$wsAms.Range("E22").Value = 119
$wsAms.Range("F22").Value = 0
$wsAms.Range("G22").Value = 2.8
$wsAms.Range("A22:G22").PasteSpecial(-4122)
